$wb = $excel.ActiveWorkbook

# ============ Sheet 1 ("展览") ============
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range('F2').Value = 362
$ws1.Range('F4').Value = 10736
$ws1.Range('F7').Value = 157
$ws1.Range('I7').Value = '//i0.hdslb.com/bfs/openplatform/202406/NlA41ysa1717751539836.png'
$ws1.Range('F8').Value = 1325
$ws1.Range('F9').Value = 8269
$ws1.Range('F11').Value = 465
$ws1.Range('F12').Value = 255
$ws1.Range('F17').Value = 327
$ws1.Range('C18').Value = '合肥·lovelive only'
$ws1.Range('D18').Value = '莲花路与繁华大道交叉口258号 百乐门大剧院'
$ws1.Range('E18').Value = '2024.07.20 11:00-07.20 19:30'
$ws1.Range('F18').Value = 12
$ws1.Range('H18').Value = 'https://show.bilibili.com/platform/detail.html?id=87058'
$ws1.Range('I18').Value = '//i2.hdslb.com/bfs/openplatform/202406/M6rzTe6y1717600620230.jpeg'
$ws1.Range('C19').Value = '合肥·第十四届次元之门动漫游戏博览会'
$ws1.Range('D19').Value = '南京路与庐州大道交汇处 合肥滨湖国际会展中心'
$ws1.Range('E19').Value = '2024.07.20 10:00-07.21 17:00'
$ws1.Range('F19').Value = 774
$ws1.Range('G19').Value = 68
$ws1.Range('H19').Value = 'https://show.bilibili.com/platform/detail.html?id=85336'
$ws1.Range('I19').Value = '//i2.hdslb.com/bfs/openplatform/202405/Bu6iQPJ01715161445356.jpeg'
$ws1.Range('B20').NumberFormat = '@'
$ws1.Range('B20').Value = '2024-07-20'
$ws1.Range('B20').Style = 'Normal'
$ws1.Range('C20').Value = '安徽·赛马娘Only 2.0'
$ws1.Range('D20').Value = '文忠路1865号 赫拉诺言艺术中心'
$ws1.Range('E20').Value = '2024.07.20 09:00-07.20 17:00'
$ws1.Range('F20').Value = 129
$ws1.Range('G20').Value = 78
$ws1.Range('H20').Value = 'https://show.bilibili.com/platform/detail.html?id=84539'
$ws1.Range('I20').Value = '//i1.hdslb.com/bfs/openplatform/202405/ibcY9Edj1715235810905.jpeg'
$ws1.Range('B21').NumberFormat = '@'
$ws1.Range('B21').Value = '2024-07-21'
$ws1.Range('B21').Style = 'Normal'
$ws1.Range('C21').Value = '合肥·首届Gumi同人展'
$ws1.Range('D21').Value = '新站区东方大道288号 少荃体育中心'
$ws1.Range('E21').Value = '2024.07.21 09:30-07.21 17:00'
$ws1.Range('F21').Value = 1064
$ws1.Range('G21').Value = 68
$ws1.Range('H21').Value = 'https://show.bilibili.com/platform/detail.html?id=86573'
$ws1.Range('I21').Value = '//i0.hdslb.com/bfs/openplatform/202405/DsvnHgmP1717038341915.jpeg'
$ws1.Range('B22').NumberFormat = '@'
$ws1.Range('B22').Value = '2024-07-27'
$ws1.Range('B22').Style = 'Normal'
$ws1.Range('C22').Value = '安徽·MAX特摄only展'
$ws1.Range('D22').Value = '桐城路127号合作经济广场3号楼23层 赤阑桥艺术空间'
$ws1.Range('E22').Value = '2024.07.27 09:30-07.27 18:00'
$ws1.Range('F22').Value = 286
$ws1.Range('G22').Value = 50
$ws1.Range('H22').Value = 'https://show.bilibili.com/platform/detail.html?id=83684'
$ws1.Range('I22').Value = '//i0.hdslb.com/bfs/openplatform/202405/qBnW1VeB1715423018997.jpeg'
$ws1.Range('B23').NumberFormat = '@'
$ws1.Range('B23').Value = '2024-07-28'
$ws1.Range('B23').Style = 'Normal'
$ws1.Range('C23').Value = '合肥·咒术回战only'
$ws1.Range('D23').Value = '清河路19号 依立腾工业园区'
$ws1.Range('E23').Value = '2024.07.28 09:30-07.28 17:30'
$ws1.Range('F23').Value = 103
$ws1.Range('G23').Value = 60
$ws1.Range('H23').Value = 'https://show.bilibili.com/platform/detail.html?id=86520'
$ws1.Range('I23').Value = '//i2.hdslb.com/bfs/openplatform/202405/cLCM0a1e1716952386781.png'
$ws1.Range('A24').Value = 23
$ws1.Range('B24').NumberFormat = '@'
$ws1.Range('B24').Value = '2024-08-03'
$ws1.Range('B24').Style = 'Normal'
$ws1.Range('C24').Value = '合肥·第七届环形宇宙动漫游戏嘉年华'
$ws1.Range('D24').Value = '南京路与庐州大道交汇处 合肥滨湖国际会展中心'
$ws1.Range('E24').Value = '2024.08.03 09:30-08.04 17:00'
$ws1.Range('F24').Value = 1756
$ws1.Range('G24').Value = 49
$ws1.Range('H24').Value = 'https://show.bilibili.com/platform/detail.html?id=84767'
$ws1.Range('I24').Value = '//i2.hdslb.com/bfs/openplatform/202404/nBGuQecO1713856894035.jpeg'

# Copy row-number cell style (bold, border, centered) onto the newly appended A24
$ws1.Range('A23').Copy()
$ws1.Range('A24').PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ============ Sheet 4 ("全部类型") ============
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range('F2').Value = 362
$ws4.Range('F4').Value = 10736
$ws4.Range('F7').Value = 157
$ws4.Range('I7').Value = '//i0.hdslb.com/bfs/openplatform/202406/NlA41ysa1717751539836.png'
$ws4.Range('F8').Value = 1325
$ws4.Range('F9').Value = 8269
$ws4.Range('F11').Value = 465
$ws4.Range('F12').Value = 255
$ws4.Range('F17').Value = 327
$ws4.Range('C18').Value = '合肥·lovelive only'
$ws4.Range('D18').Value = '莲花路与繁华大道交叉口258号 百乐门大剧院'
$ws4.Range('E18').Value = '2024.07.20 11:00-07.20 19:30'
$ws4.Range('F18').Value = 12
$ws4.Range('H18').Value = 'https://show.bilibili.com/platform/detail.html?id=87058'
$ws4.Range('I18').Value = '//i2.hdslb.com/bfs/openplatform/202406/M6rzTe6y1717600620230.jpeg'
$ws4.Range('C19').Value = '合肥·第十四届次元之门动漫游戏博览会'
$ws4.Range('D19').Value = '南京路与庐州大道交汇处 合肥滨湖国际会展中心'
$ws4.Range('E19').Value = '2024.07.20 10:00-07.21 17:00'
$ws4.Range('F19').Value = 774
$ws4.Range('G19').Value = 68
$ws4.Range('H19').Value = 'https://show.bilibili.com/platform/detail.html?id=85336'
$ws4.Range('I19').Value = '//i2.hdslb.com/bfs/openplatform/202405/Bu6iQPJ01715161445356.jpeg'
$ws4.Range('B20').NumberFormat = '@'
$ws4.Range('B20').Value = '2024-07-20'
$ws4.Range('B20').Style = 'Normal'
$ws4.Range('C20').Value = '安徽·赛马娘Only 2.0'
$ws4.Range('D20').Value = '文忠路1865号 赫拉诺言艺术中心'
$ws4.Range('E20').Value = '2024.07.20 09:00-07.20 17:00'
$ws4.Range('F20').Value = 129
$ws4.Range('G20').Value = 78
$ws4.Range('H20').Value = 'https://show.bilibili.com/platform/detail.html?id=84539'
$ws4.Range('I20').Value = '//i1.hdslb.com/bfs/openplatform/202405/ibcY9Edj1715235810905.jpeg'
$ws4.Range('B21').NumberFormat = '@'
$ws4.Range('B21').Value = '2024-07-21'
$ws4.Range('B21').Style = 'Normal'
$ws4.Range('C21').Value = '合肥·首届Gumi同人展'
$ws4.Range('D21').Value = '新站区东方大道288号 少荃体育中心'
$ws4.Range('E21').Value = '2024.07.21 09:30-07.21 17:00'
$ws4.Range('F21').Value = 1064
$ws4.Range('G21').Value = 68
$ws4.Range('H21').Value = 'https://show.bilibili.com/platform/detail.html?id=86573'
$ws4.Range('I21').Value = '//i0.hdslb.com/bfs/openplatform/202405/DsvnHgmP1717038341915.jpeg'
$ws4.Range('B22').NumberFormat = '@'
$ws4.Range('B22').Value = '2024-07-27'
$ws4.Range('B22').Style = 'Normal'
$ws4.Range('C22').Value = '安徽·MAX特摄only展'
$ws4.Range('D22').Value = '桐城路127号合作经济广场3号楼23层 赤阑桥艺术空间'
$ws4.Range('E22').Value = '2024.07.27 09:30-07.27 18:00'
$ws4.Range('F22').Value = 286
$ws4.Range('G22').Value = 50
$ws4.Range('H22').Value = 'https://show.bilibili.com/platform/detail.html?id=83684'
$ws4.Range('I22').Value = '//i0.hdslb.com/bfs/openplatform/202405/qBnW1VeB1715423018997.jpeg'
$ws4.Range('B23').NumberFormat = '@'
$ws4.Range('B23').Value = '2024-07-28'
$ws4.Range('B23').Style = 'Normal'
$ws4.Range('C23').Value = '合肥·咒术回战only'
$ws4.Range('D23').Value = '清河路19号 依立腾工业园区'
$ws4.Range('E23').Value = '2024.07.28 09:30-07.28 17:30'
$ws4.Range('F23').Value = 103
$ws4.Range('G23').Value = 60
$ws4.Range('H23').Value = 'https://show.bilibili.com/platform/detail.html?id=86520'
$ws4.Range('I23').Value = '//i2.hdslb.com/bfs/openplatform/202405/cLCM0a1e1716952386781.png'
$ws4.Range('C24').Value = '合肥·第七届环形宇宙动漫游戏嘉年华'
$ws4.Range('D24').Value = '南京路与庐州大道交汇处 合肥滨湖国际会展中心'
$ws4.Range('E24').Value = '2024.08.03 09:30-08.04 17:00'
$ws4.Range('F24').Value = 1756
$ws4.Range('G24').Value = 49
$ws4.Range('H24').Value = 'https://show.bilibili.com/platform/detail.html?id=84767'
$ws4.Range('I24').Value = '//i2.hdslb.com/bfs/openplatform/202404/nBGuQecO1713856894035.jpeg'
$ws4.Range('A25').Value = 24
$ws4.Range('B25').NumberFormat = '@'
$ws4.Range('B25').Value = '2024-08-03'
$ws4.Range('B25').Style = 'Normal'
$ws4.Range('C25').Value = '合肥·首届包河留声机音乐节—《菊次郎的夏天》久石让钢琴曲梦幻之旅演奏会'
$ws4.Range('D25').Value = '徽州大道辅路与祁门路辅路交叉口北120米 包河凤凰剧院'
$ws4.Range('E25').Value = '2024.08.03 19:30-08.03 21:00'
$ws4.Range('F25').Value = 30
$ws4.Range('G25').Value = 80
$ws4.Range('H25').Value = 'https://show.bilibili.com/platform/detail.html?id=83556'
$ws4.Range('I25').Value = '//i1.hdslb.com/bfs/openplatform/202403/4nwOTVDu1711695345941.jpeg'

# Copy row-number cell style (bold, border, centered) onto the newly appended A25
$ws4.Range('A24').Copy()
$ws4.Range('A25').PasteSpecial(-4122)
$excel.CutCopyMode = $false

